# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# with latest scraped values (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.684.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "'3.798.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'596.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "'167.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.24%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "'0.0000253"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "'4.440.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").Value = "'3.781.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").Value = "'18.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.43%  "
$ws.Range("D17").Value = "'67.695.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("E18").Value = "  +2.97%  "
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").Value = "'461.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").Value = "'9.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.97%  "
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("E23").Value = "  +3.99%  "
$ws.Range("D24").Value = "'83.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").Value = "'3.941.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("E31").Value = "  +3.19%  "
$ws.Range("D32").Value = "'7.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("D33").Value = "'29.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "'9.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").Value = "'0.0999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("E37").Value = "  +3.29%  "
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("D40").Value = "'5.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D43").Value = "'48.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.31%  "
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("D45").Value = "'43.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "'148.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").Value = "'395.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("D50").Value = "'26.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.68%  "
$ws.Range("D51").Value = "'1.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.12%  "
